$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Most of these replacement values parse as plain numbers (e.g. "570.44"), but the
# source data stores Price as literal text (note values like "64.447.56" which are not
# valid numbers at all). Force the cells to Text format before assigning so Excel keeps
# them as strings, then clear the temporary number format again so no extra cell style
# is left behind (the workbook ships these cells with the default/no style).
$priceCells = [ordered]@{
    "D2" = "64.447.56"
    "D3" = "3.165.21"
    "D5" = "570.44"
    "D6" = "164.04"
    "D8" = "0.585"
    "D10" = "6.64"
    "D11" = "0.385"
    "D12" = "3.714.42"
    "D14" = "64.517.75"
    "D15" = "25.32"
    "D16" = "3.160.14"
    "D17" = "0.0000156"
    "D18" = "407.71"
    "D19" = "12.74"
    "D20" = "5.29"
    "D21" = "7.10"
    "D23" = "68.56"
    "D24" = "0.484"
    "D25" = "0.197"
    "D26" = "0.0000103"
    "D27" = "8.91"
    "D28" = "1.00"
    "D29" = "1.81"
    "D30" = "21.21"
    "D31" = "6.36"
    "D32" = "4.91"
    "D33" = "1.13"
    "D34" = "155.56"
    "D35" = "1.35"
    "D36" = "2.691.83"
    "D37" = "1.69"
    "D38" = "24.15"
    "D39" = "4.09"
    "D40" = "0.695"
    "D41" = "0.0621"
    "D42" = "5.46"
    "D43" = "0.0259"
    "D44" = "292.11"
    "D45" = "21.53"
    "D47" = "0.0985"
    "D48" = "1.93"
    "D49" = "10.51"
    "D50" = "5.73"
    "D51" = "0.885"
}
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"
foreach ($addr in $priceCells.Keys) {
    $ws.Range($addr).Value = $priceCells[$addr]
}
$priceRange.ClearFormats()

# --- Coin / Link swap for rows 43-45 (VeChain, Bittensor, InjectiveProtocol re-ranked) ---
$coinCells = [ordered]@{
    "B43" = "VeChain"
    "B44" = "Bittensor"
    "B45" = "InjectiveProtocol"
    "C43" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "C44" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "C45" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
}
foreach ($addr in $coinCells.Keys) {
    $ws.Range($addr).Value = $coinCells[$addr]
}

# --- Column E (Volume(1h)) updates ---
$volumeCells = [ordered]@{
    "E2" = "  +0.73%  "
    "E3" = "  -0.72%  "
    "E4" = "  +0.04%  "
    "E5" = "  +0.01%  "
    "E6" = "  -3.27%  "
    "E7" = "  +0.03%  "
    "E8" = "  -4.62%  "
    "E9" = "  -3.13%  "
    "E10" = "  -1.35%  "
    "E11" = "  -0.57%  "
    "E12" = "  -0.59%  "
    "E13" = "  -1.10%  "
    "E14" = "  +0.66%  "
    "E15" = "  -0.57%  "
    "E16" = "  -1.09%  "
    "E17" = "  -2.28%  "
    "E18" = "  -2.00%  "
    "E19" = "  -0.87%  "
    "E20" = "  -1.50%  "
    "E21" = "  -1.18%  "
    "E22" = "  +0.31%  "
    "E23" = "  -2.96%  "
    "E24" = "  -1.83%  "
    "E25" = "  -2.86%  "
    "E26" = "  -6.22%  "
    "E27" = "  +1.27%  "
    "E28" = "  +0.02%  "
    "E29" = "  -1.90%  "
    "E30" = "  -3.32%  "
    "E31" = "  -0.82%  "
    "E32" = "  -2.36%  "
    "E33" = "  -0.96%  "
    "E34" = "  -0.47%  "
    "E35" = "  -2.04%  "
    "E36" = "  -1.57%  "
    "E37" = "  -1.05%  "
    "E38" = "  -4.03%  "
    "E39" = "  -2.26%  "
    "E40" = "  -3.17%  "
    "E41" = "  -1.42%  "
    "E42" = "  -4.48%  "
    "E43" = "  -1.59%  "
    "E44" = "  -1.93%  "
    "E45" = "  -2.83%  "
    "E46" = "  +0.02%  "
    "E47" = "  -1.08%  "
    "E48" = "  -8.07%  "
    "E49" = "  +0.91%  "
    "E50" = "  -1.55%  "
    "E51" = "  -4.91%  "
}
foreach ($addr in $volumeCells.Keys) {
    $ws.Range($addr).Value = $volumeCells[$addr]
}
